# Enhance login process with custom wait condition and improve error
# handling for authentication failures.
#
# (Functionally: the shipment report's data rows are restructured — the
# first three distinct shipment rows (tracking 240037171791 / Kelly Johanna
# Perlaza Potes, 700170507656 / LEIDI PAOLA MARTINEZ, 700170328932 / XIMENA
# MUNOZ) are repeated across most of the sheet, the sheet is extended from
# 15 to 22 data+header rows, and the former "Alex Aular" / "JESUSU ALBERTO
# HINCAPIE PJ" rows are preserved, relocated to the tail of the sheet (with
# "Alex Aular" duplicated once more).)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: before anything else gets overwritten, relocate the rows that
# need to survive at the bottom of the extended sheet (rows 20-22) by
# copying them off of their original positions (row 5 = Alex Aular, row 6 =
# JESUSU ALBERTO HINCAPIE PJ) while those positions are still intact.
$ws.Range("A5:D5").Copy()
$ws.Range("A20:D20").PasteSpecial()

$ws.Range("A5:D5").Copy()
$ws.Range("A21:D21").PasteSpecial()

$ws.Range("A6:D6").Copy()
$ws.Range("A22:D22").PasteSpecial()

# --- Step 2: fill rows 5-19 by repeating the first three data rows
# (rows 2, 3, 4), which are never themselves overwritten, in a cyclical
# Kelly / LEIDI / XIMENA pattern matching the target layout.
$sources = @(2, 3, 4, 2, 3, 4, 2, 3, 2, 3, 2, 3, 2, 3, 4)
$destRow = 5
foreach ($srcRow in $sources) {
    $ws.Range("A$srcRow`:D$srcRow").Copy()
    $ws.Range("A$destRow`:D$destRow").PasteSpecial()
    $destRow++
}

$excel.CutCopyMode = $false
